$d = $word.ActiveDocument

# Map of old text -> new text (unique strings in the document, so a
# whole-word, case-sensitive, non-wildcard Find/Replace is safe).
$replacements = [ordered]@{
    "2024-09-06 Friday" = "2024-09-07 Saturday"
    "31×74="            = "89×70="
    "59×46="            = "60×79="
    "65×93="            = "68×21="
    "37×93="            = "77×43="
    "43×64="            = "78×96="
    "13×63="            = "15×57="
    "47×91="            = "79×34="
    "55×41="            = "26×41="
    "15×21="            = "19×18="
    "39×60="            = "82×17="
    "81×49="            = "20×73="
    "34×97="            = "70×36="
    "43×84="            = "43×41="
    "45×18="            = "19×35="
    "81×43="            = "49×78="
    "17×49="            = "21×44="
    "45×23="            = "64×24="
    "45×47="            = "86×44="
    "96×77="            = "44×51="
    "69×85="            = "88×73="
    "97×29="            = "36×97="
    "71×62="            = "14×25="
    "93×41="            = "98×82="
    "12×85="            = "90×84="
    "83×59="            = "92×94="
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
